$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 42 and 43: Coin (B) and Link (C) swapped between RenderToken and Algorand ---
$ws.Range("B42").Value = "Algorand"
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"

# --- Price (column D) updates. ---
# These are stored as plain text in the workbook (e.g. "29.667.73", "0.00001028")
# A leading apostrophe forces Excel to keep the literal text instead of
# silently re-parsing it as a number (which would mangle the precision/
# formatting, e.g. "72.00" -> 72). Style is reset to Normal afterwards so no
# stray quote-prefix/number-format style gets attached to the cell.
$ws.Range("D2").Value = "'29.681.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.924.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Value = "'339.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = "'0.4820"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.4066"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.08109"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'1.003"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'23.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'1.963.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'5.996"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'7.200"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'90.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.06858"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'1.013"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.00001030"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'17.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'29.692.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'5.562"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'11.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.163"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'2.199.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'6.604"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'157.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'19.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'2.069"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'120.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'1.004"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'0.09603"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'5.534"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'1.401"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Value = "'0.06535"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'0.02259"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'1.199"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'0.5907"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'10.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'7.867"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1834"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'2.478"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'1.245"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'12.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.07473"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.5527"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'1.965"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'116.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'2.404"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'72.00"
$ws.Range("D51").Style = "Normal"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.64%  "
$ws.Range("E5").Value = "  +4.49%  "
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("E22").Value = "  -1.60%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("E29").Value = "  -2.09%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +7.01%  "
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("E44").Value = "  -2.65%  "
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("E51").Value = "  -0.62%  "
